# Updated cryptos list on Thu Aug 24 10:11:39 UTC 2023 with GitHub Actions
# Applies the per-cell price/volume(1h) refresh described by the diff.
# Two coin rows (12/13 and 46/47) were also re-ordered/swapped upstream,
# so those rows' Coin/Link/Price/Volume cells are fully rewritten too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '26.495.05'
$ws.Range("E2").Value = '  +1.76%  '

# Row 3
$ws.Range("E3").Value = '  +1.56%  '

# Row 4
$ws.Range("E4").Value = '  +0.03%  '

# Row 5
$ws.Range("D5").Value = '219.94'
$ws.Range("E5").Value = '  +2.06%  '

# Row 6
$ws.Range("D6").Value = '0.5270'
$ws.Range("E6").Value = '  +0.94%  '

# Row 7
$ws.Range("E7").Value = '  +0.03%  '

# Row 8
$ws.Range("D8").Value = '0.2676'
$ws.Range("E8").Value = '  +2.52%  '

# Row 9
$ws.Range("D9").Value = '0.06372'
$ws.Range("E9").Value = '  +0.21%  '

# Row 10
$ws.Range("D10").Value = '21.72'
$ws.Range("E10").Value = '  +4.30%  '

# Row 11
$ws.Range("D11").Value = '0.07796'
$ws.Range("E11").Value = '  +1.76%  '

# Row 12
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = '1.674.48'
$ws.Range("E12").Value = '  +1.71%  '

# Row 13
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").Value = '4.486'
$ws.Range("E13").Value = '  +1.47%  '

# Row 14
$ws.Range("D14").Value = '0.5557'
$ws.Range("E14").Value = '  +0.25%  '

# Row 15
$ws.Range("D15").Value = '0.0₅8276'
$ws.Range("E15").Value = '  -0.38%  '

# Row 16
$ws.Range("D16").Value = '65.59'
$ws.Range("E16").Value = '  +1.09%  '

# Row 17
$ws.Range("D17").Value = '26.488.22'
$ws.Range("E17").Value = '  +1.67%  '

# Row 18
$ws.Range("E18").Value = '  +0.03%  '

# Row 19
$ws.Range("D19").Value = '4.755'
$ws.Range("E19").Value = '  +0.66%  '

# Row 20
$ws.Range("D20").Value = '192.80'
$ws.Range("E20").Value = '  +2.35%  '

# Row 21
$ws.Range("D21").Value = '10.32'
$ws.Range("E21").Value = '  +1.32%  '

# Row 22
$ws.Range("D22").Value = '6.294'
$ws.Range("E22").Value = '  +0.56%  '

# Row 23
$ws.Range("E23").Value = '  +0.10%  '

# Row 24
$ws.Range("E24").Value = '  +4.06%  '

# Row 25
$ws.Range("D25").Value = '137.94'
$ws.Range("E25").Value = '  -5.40%  '

# Row 26
$ws.Range("D26").Value = '7.378'
$ws.Range("E26").Value = '  -0.46%  '

# Row 27
$ws.Range("E27").Value = '  +2.92%  '

# Row 28
$ws.Range("D28").Value = '1.430'
$ws.Range("E28").Value = '  +2.08%  '

# Row 29
$ws.Range("E29").Value = '  +4.72%  '

# Row 30
$ws.Range("D30").Value = '1.289'
$ws.Range("E30").Value = '  +1.82%  '

# Row 31
$ws.Range("D31").Value = '3.614'
$ws.Range("E31").Value = '  +6.06%  '

# Row 32
$ws.Range("D32").Value = '3.415'
$ws.Range("E32").Value = '  +0.52%  '

# Row 33
$ws.Range("D33").Value = '1.685'
$ws.Range("E33").Value = '  +1.86%  '

# Row 34
$ws.Range("D34").Value = '1.006'
$ws.Range("E34").Value = '  +0.92%  '

# Row 35
$ws.Range("D35").Value = '0.6147'
$ws.Range("E35").Value = '  +9.25%  '

# Row 36
$ws.Range("D36").Value = '2.423'
$ws.Range("E36").Value = '  +1.24%  '

# Row 37
$ws.Range("D37").Value = '2.782'
$ws.Range("E37").Value = '  +1.03%  '

# Row 38
$ws.Range("D38").Value = '0.01612'
$ws.Range("E38").Value = '  +0.17%  '

# Row 39
$ws.Range("D39").Value = '6.023'
$ws.Range("E39").Value = '  +2.90%  '

# Row 40
$ws.Range("D40").Value = '1.092.63'
$ws.Range("E40").Value = '  +6.38%  '

# Row 41
$ws.Range("D41").Value = '0.8594'
$ws.Range("E41").Value = '  +0.45%  '

# Row 42
$ws.Range("E42").Value = '  -0.02%  '

# Row 43
$ws.Range("D43").Value = '100.68'
$ws.Range("E43").Value = '  +2.15%  '

# Row 44
$ws.Range("D44").Value = '1.815.41'
$ws.Range("E44").Value = '  +1.10%  '

# Row 45
$ws.Range("D45").Value = '58.53'
$ws.Range("E45").Value = '  +5.05%  '

# Row 46
$ws.Range("B46").Value = 'BabyDogeCoin'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D46").Value = '0.0₈106'
$ws.Range("E46").Value = '  -5.35%  '

# Row 47
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").Value = '8.140'
$ws.Range("E47").Value = '  +0.57%  '

# Row 48
$ws.Range("D48").Value = '0.9989'
$ws.Range("E48").Value = '  -0.59%  '

# Row 49
$ws.Range("D49").Value = '1.514'
$ws.Range("E49").Value = '  +9.53%  '

# Row 50
$ws.Range("D50").Value = '0.05192'
$ws.Range("E50").Value = '  +0.81%  '

# Row 51
$ws.Range("E51").Value = '  +0.39%  '

